$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column I ("roboticS1Prep") currently holds the text "No" for every
# data row (2-39). Convert those to a real boolean FALSE, formatted so
# it still displays as text ("TRUE"/"FALSE").

# Rows 2-27 already share one explicit style; rows 28-39 use the sheet's
# default (unstyled) cells. Unify every I2:I39 cell onto the same visual
# style first (copy format from I2), so the whole column collapses onto a
# single cell style after the number-format change below.
$ws.Range("I2").Copy() | Out-Null
$ws.Range("I28:I39").PasteSpecial(-4122) | Out-Null

$iCol = $ws.Range("I2:I39")
$iCol.Value = $false
$iCol.NumberFormat = '"TRUE";"TRUE";"FALSE"'

# Match the author's new selection/scroll position recorded in the sheet.
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I2:I39").Select() | Out-Null
